$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($ws, $row, $values)
    $n = $values.Count
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $c1 = $ws.Cells.Item($row, 1)
    $c2 = $ws.Cells.Item($row, $n)
    $rng = $ws.Range($c1, $c2)
    $rng.Value = $arr
}

# Weekly refresh: two new weekly records are inserted right after the
# existing row 85 ("Feria Lagunitas de Puerto Montt" / Kiwi series),
# pushing the rest of the historical rows (old 88..111) down by two rows.
$ws.Rows("86:87").Insert()

Set-RowValues $ws 86 @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44463, 10, "Fruta", 100101, "Berries", 100101007, "Kiwi", "Hayward", "Especial", 300, 20000, 20000, 20000, "`$/caja 15 kilos", "Provincia de Curicó", 1333, 15)
Set-RowValues $ws 87 @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44463, 10, "Fruta", 100101, "Berries", 100101007, "Kiwi", "Hayward", "Primera", 150, 15000, 15000, 15000, "`$/caja 15 kilos", "Provincia de Curicó", 1000, 15)
